# Adding commentary implementations for other charts. Resolves #60.
#
# - Drop Sheet3 entirely, repurpose Sheet2 as a "Commentary" sheet with a
#   small region/topic commentary table.
# - The deletion of Sheet3 shifts the active sheet to (what becomes)
#   "Commentary", matching the workbook's new activeTab.

$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("Sheet3").Delete() | Out-Null

$ws = $wb.Worksheets.Item("Sheet2")
$ws.Name = "Commentary"

# Header row
$ws.Range("A1").Value = "Region"
$ws.Range("B1").Value = "Cane and Horticulture"
$ws.Range("C1").Value = "Grazing"
$ws.Range("D1").Value = "Nitrogen and Pesticide"
$ws.Range("E1").Value = "Sediment"

# Data row
$ws.Range("A2").Value = "GBR"
$ws.Range("A2").VerticalAlignment = -4160

$ws.Range("B2:E2").Value = "This is the **Cane and Horticulture** commentary."
$ws.Range("B2:E2").VerticalAlignment = -4160
$ws.Range("B2:E2").WrapText = $true

$ws.Range("C2").Value = "_Nothing yet._"
$ws.Range("D2").Value = "_Nothing yet._"
$ws.Range("E2").Value = "_Nothing yet._"

$ws.Rows.Item(2).RowHeight = 340.25

$ws.Columns.Item(1).ColumnWidth = 7.22
$ws.Range("B1:E1").ColumnWidth = 35.72

$ws.Range("C2").Select() | Out-Null
